{"js": "// Apply the benchmark-stats update to the single-column results table.\n// Each row holds one stat value in its only cell; a handful of rows are\n// being corrected, and the last three \"raw dump\" rows (which currently\n// hold a whole tab-separated record) are being collapsed down to the\n// single summary value that used to live in the first three rows.\n\nconst table = context.document.body.tables.getFirst();\ntable.rows.load(\"items\");\nawait context.sync();\n\nconst rows = table.rows.items;\nfor (const row of rows) {\n  row.cells.load(\"items\");\n}\nawait context.sync();\n\n// index -> new cell value (only rows that actually change are listed)\nconst updates = {\n  0: \"0M\",\n  1: \"0M\",\n  2: \"0M\",\n  3: \"804\",\n  5: \"0.00060\",\n  6: \"0.00015\",\n  7: \"0.00004\",\n  8: \"0.00021\",\n  9: \"0.00022\",\n  10: \"0.00035\",\n  11: \"0.11676\",\n  43: \"99.97\",\n  44: \"0.12\",\n  45: \"383\",\n};\n\nfor (const [idxStr, newValue] of Object.entries(updates)) {\n  const idx = Number(idxStr);\n  const cell = rows[idx].cells.items[0];\n  cell.value = newValue;\n}\n\nawait context.sync();\n", "ps1": "# Apply the benchmark-stats update to the single-column results table.\n# Each row holds one stat value in its only cell; a handful of rows are\n# being corrected, and the last three \"raw dump\" rows (which currently\n# hold a whole tab-separated record) are being collapsed down to the\n# single summary value that used to live in the first three rows.\n\n$doc = $word.ActiveDocument\n$table = $doc.Tables.Item(1)\n\n# 1-based row index -> new cell value (only rows that actually change)\n$updates = [ordered]@{\n    1  = \"0M\"\n    2  = \"0M\"\n    3  = \"0M\"\n    4  = \"804\"\n    6  = \"0.00060\"\n    7  = \"0.00015\"\n    8  = \"0.00004\"\n    9  = \"0.00021\"\n    10 = \"0.00022\"\n    11 = \"0.00035\"\n    12 = \"0.11676\"\n    44 = \"99.97\"\n    45 = \"0.12\"\n    46 = \"383\"\n}\n\nforeach ($rowIndex in $updates.Keys) {\n    $cell = $table.Cell($rowIndex, 1)\n    $cell.Range.Text = $updates[$rowIndex]\n}\n"}
